$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("acc_final")
$ws.Range("A1").Value = "test"
